# Refresh the crypto "Price" (D) / "Volume(1h)" (E) columns on Sheet1 to the
# latest scraped values (coinranking.com snapshot), matching the GitHub
# Actions cron-job commit "Updated symbol list".
#
# The sheet stores these figures as literal text (e.g. "5.063", "1.13%"),
# not numbers, so each write is forced to Text via a leading apostrophe and
# the cell's style is immediately reset to "Normal" afterwards so that no
# stray number-format / quote-prefix styling is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "E3";  Value = "1.32%" },
    @{ Cell = "D4";  Value = "5.063" },
    @{ Cell = "E4";  Value = "1.13%" },
    @{ Cell = "D5";  Value = "0.08118" },
    @{ Cell = "E5";  Value = "0.48%" },
    @{ Cell = "D6";  Value = "2.038" },
    @{ Cell = "E6";  Value = "5.63%" },
    @{ Cell = "D7";  Value = "4.162" },
    @{ Cell = "D8";  Value = "7.865" },
    @{ Cell = "D9";  Value = "0.9268" },
    @{ Cell = "E9";  Value = "-0.48%" },
    @{ Cell = "D10"; Value = "0.1429" },
    @{ Cell = "E10"; Value = "14.31%" },
    @{ Cell = "D11"; Value = "0.1924" },
    @{ Cell = "E11"; Value = "0.54%" },
    @{ Cell = "D12"; Value = "0.09135" },
    @{ Cell = "E12"; Value = "-0.80%" },
    @{ Cell = "E13"; Value = "-1.97%" },
    @{ Cell = "D14"; Value = "0.09903" },
    @{ Cell = "E14"; Value = "-0.29%" },
    @{ Cell = "D15"; Value = "0.001404" },
    @{ Cell = "E15"; Value = "-0.75%" },
    @{ Cell = "D16"; Value = "0.006159" },
    @{ Cell = "E16"; Value = "-8.45%" },
    @{ Cell = "D17"; Value = "3.838" },
    @{ Cell = "E17"; Value = "6.09%" },
    @{ Cell = "E18"; Value = "3.65%" },
    @{ Cell = "D19"; Value = "0.3440" },
    @{ Cell = "E19"; Value = "-0.05%" },
    @{ Cell = "E21"; Value = "-7.12%" },
    @{ Cell = "E22"; Value = "-7.55%" },
    @{ Cell = "D23"; Value = "0.04367" },
    @{ Cell = "E23"; Value = "-0.94%" },
    @{ Cell = "D24"; Value = "0.001232" },
    @{ Cell = "E24"; Value = "-0.22%" },
    @{ Cell = "E25"; Value = "4.24%" },
    @{ Cell = "D27"; Value = "0.0001300" },
    @{ Cell = "E27"; Value = "-0.09%" },
    @{ Cell = "D39"; Value = "0.02035" },
    @{ Cell = "E39"; Value = "3.57%" },
    @{ Cell = "D40"; Value = "0.05158" },
    @{ Cell = "E40"; Value = "-0.30%" },
    @{ Cell = "D41"; Value = "0.007491" },
    @{ Cell = "E41"; Value = "-1.18%" },
    @{ Cell = "D42"; Value = "0.01011" },
    @{ Cell = "E42"; Value = "-0.33%" },
    @{ Cell = "D43"; Value = "0.1373" },
    @{ Cell = "E43"; Value = "0.25%" },
    @{ Cell = "D44"; Value = "0.002130" },
    @{ Cell = "E44"; Value = "1.34%" },
    @{ Cell = "E45"; Value = "-8.87%" },
    @{ Cell = "D46"; Value = "0.00006291" },
    @{ Cell = "E46"; Value = "-1.36%" },
    @{ Cell = "E47"; Value = "-0.09%" },
    @{ Cell = "D48"; Value = "64.86" },
    @{ Cell = "E48"; Value = "-0.16%" },
    @{ Cell = "E49"; Value = "-22.04%" },
    @{ Cell = "E50"; Value = "-0.09%" },
    @{ Cell = "E51"; Value = "-0.09%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}
